# "Status code check 2.0"
# - Adds a new "statusCode" column (between "errorMessage" and "timestamp")
#   to the surviving sheets, populated with either an HTTP status code
#   (as a number) or the literal text "Unknown".
# - Refreshes the "timestamp" column with the new run's values.
# - Removes the "LMBC" and "Abbvie Pro Medical" sheets, leaving
#   Sitegen / Allerganpro / ADPA (in that order).

$wb = $excel.ActiveWorkbook

# --- Sitegen (Healthcareheroes) ---------------------------------------
$ws = $wb.Worksheets("Sitegen")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "statusCode"
$ws.Range("D2").Value = "Unknown"
$ws.Range("E2").Value = "2024-06-21T06:14:14.947Z"

# --- Allerganpro (India) -----------------------------------------------
$ws = $wb.Worksheets("Allerganpro")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "statusCode"
$ws.Range("D2").Value = 404
$ws.Range("E2").Value = "2024-06-21T06:17:08.073Z"

# --- ADPA (Belgium / Paraguay / Uruguay) --------------------------------
$ws = $wb.Worksheets("ADPA")
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "statusCode"
$ws.Range("D2").Value = "Unknown"
$ws.Range("E2").Value = "2024-06-21T06:20:25.392Z"
$ws.Range("D3").Value = 404
$ws.Range("E3").Value = "2024-06-21T06:20:42.174Z"
$ws.Range("D4").Value = 404
$ws.Range("E4").Value = "2024-06-21T06:20:58.782Z"

# --- Drop the retired sheets --------------------------------------------
$excel.DisplayAlerts = $false
[void]$wb.Worksheets("LMBC").Delete()
[void]$wb.Worksheets("Abbvie Pro Medical").Delete()
